$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D for "Status", shifting Jan_2026..QoQ right by one
$ws.Columns("D:D").Insert()

# Insert 4 new rows before row 28 to make room for the new holdings rows
$ws.Rows("28:31").Insert()

# Header row
$ws.Cells.Item(1,4).Value = "Status"
$ws.Cells.Item(1,5).Value = "Jan_2026"
$ws.Cells.Item(1,6).Value = "Dec_2025"
$ws.Cells.Item(1,7).Value = "Oct_2025"
$ws.Cells.Item(1,8).Value = "MoM"
$ws.Cells.Item(1,9).Value = "QoQ"

# Row 2: INE090A01021 - ICICI Bank Limited
$ws.Cells.Item(2,1).Value = "INE090A01021"
$ws.Cells.Item(2,2).Value = "ICICI Bank Limited"
$ws.Cells.Item(2,3).Value = "quant Momentum Fund"
$ws.Cells.Item(2,4).Value = "Fresh Entry"
$ws.Cells.Item(2,5).Value = 9.564712
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(2,7).Value = 0
$ws.Cells.Item(2,8).Value = 9.564712
$ws.Cells.Item(2,9).Value = 9.564712

# Row 3: INE040A01034 - HDFC Bank Limited
$ws.Cells.Item(3,1).Value = "INE040A01034"
$ws.Cells.Item(3,2).Value = "HDFC Bank Limited"
$ws.Cells.Item(3,3).Value = "quant Momentum Fund"
$ws.Cells.Item(3,4).Value = "Fresh Entry"
$ws.Cells.Item(3,5).Value = 8.956498
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(3,7).Value = 0
$ws.Cells.Item(3,8).Value = 8.956498
$ws.Cells.Item(3,9).Value = 8.956498

# Row 4: INE423A01024 - Adani Enterprises Limited
$ws.Cells.Item(4,1).Value = "INE423A01024"
$ws.Cells.Item(4,2).Value = "Adani Enterprises Limited"
$ws.Cells.Item(4,3).Value = "quant Momentum Fund"
$ws.Cells.Item(4,4).Value = "Reducing"
$ws.Cells.Item(4,5).Value = 8.873133
$ws.Cells.Item(4,6).Value = 9.170431
$ws.Cells.Item(4,7).Value = 4.50724
$ws.Cells.Item(4,8).Value = -0.2972980000000014
$ws.Cells.Item(4,9).Value = 4.365892999999999

# Row 5: INE202B01038 - Piramal Finance Ltd
$ws.Cells.Item(5,1).Value = "INE202B01038"
$ws.Cells.Item(5,2).Value = "Piramal Finance Ltd"
$ws.Cells.Item(5,3).Value = "quant Momentum Fund"
$ws.Cells.Item(5,4).Value = "Reducing Consistently"
$ws.Cells.Item(5,5).Value = 8.520021
$ws.Cells.Item(5,6).Value = 9.898523
$ws.Cells.Item(5,7).Value = 9.050035
$ws.Cells.Item(5,8).Value = -1.378502000000001
$ws.Cells.Item(5,9).Value = -0.5300139999999995

# Row 6: INE406A01037 - Aurobindo Pharma Limited
$ws.Cells.Item(6,1).Value = "INE406A01037"
$ws.Cells.Item(6,2).Value = "Aurobindo Pharma Limited"
$ws.Cells.Item(6,3).Value = "quant Momentum Fund"
$ws.Cells.Item(6,4).Value = "Adding Consistently"
$ws.Cells.Item(6,5).Value = 6.533598
$ws.Cells.Item(6,6).Value = 5.966757
$ws.Cells.Item(6,7).Value = 5.470548
$ws.Cells.Item(6,8).Value = 0.5668409999999993
$ws.Cells.Item(6,9).Value = 1.06305

# Row 7: INE364U01010 - Adani Green Energy Limited
$ws.Cells.Item(7,1).Value = "INE364U01010"
$ws.Cells.Item(7,2).Value = "Adani Green Energy Limited"
$ws.Cells.Item(7,3).Value = "quant Momentum Fund"
$ws.Cells.Item(7,4).Value = "Reducing Consistently"
$ws.Cells.Item(7,5).Value = 5.842646
$ws.Cells.Item(7,6).Value = 6.487632
$ws.Cells.Item(7,7).Value = 8.581283
$ws.Cells.Item(7,8).Value = -0.6449859999999994
$ws.Cells.Item(7,9).Value = -2.738637000000001

# Row 8: INE795G01014 - HDFC Life Insurance Co Ltd
$ws.Cells.Item(8,1).Value = "INE795G01014"
$ws.Cells.Item(8,2).Value = "HDFC Life Insurance Co Ltd"
$ws.Cells.Item(8,3).Value = "quant Momentum Fund"
$ws.Cells.Item(8,4).Value = "Adding Consistently"
$ws.Cells.Item(8,5).Value = 5.796873
$ws.Cells.Item(8,6).Value = 5.543845
$ws.Cells.Item(8,7).Value = 5.152885
$ws.Cells.Item(8,8).Value = 0.2530279999999996
$ws.Cells.Item(8,9).Value = 0.6439879999999993

# Row 9: INE917I01010 - Bajaj Auto Limited
$ws.Cells.Item(9,1).Value = "INE917I01010"
$ws.Cells.Item(9,2).Value = "Bajaj Auto Limited"
$ws.Cells.Item(9,3).Value = "quant Momentum Fund"
$ws.Cells.Item(9,4).Value = "Adding Consistently"
$ws.Cells.Item(9,5).Value = 3.971147
$ws.Cells.Item(9,6).Value = 3.604163
$ws.Cells.Item(9,7).Value = 3.266884
$ws.Cells.Item(9,8).Value = 0.3669840000000004
$ws.Cells.Item(9,9).Value = 0.7042630000000001

# Row 10: INE237A01036 - Kotak Mahindra Bank Limited
$ws.Cells.Item(10,1).Value = "INE237A01036"
$ws.Cells.Item(10,2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(10,3).Value = "quant Momentum Fund"
$ws.Cells.Item(10,4).Value = "Fresh Entry"
$ws.Cells.Item(10,5).Value = 3.67925
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(10,7).Value = 0
$ws.Cells.Item(10,8).Value = 3.67925
$ws.Cells.Item(10,9).Value = 3.67925

# Row 11: INE200M01039 - Varun Beverages Limited
$ws.Cells.Item(11,1).Value = "INE200M01039"
$ws.Cells.Item(11,2).Value = "Varun Beverages Limited"
$ws.Cells.Item(11,3).Value = "quant Momentum Fund"
$ws.Cells.Item(11,4).Value = "Fresh Entry"
$ws.Cells.Item(11,5).Value = 2.929109
$ws.Cells.Item(11,6).Value = 0
$ws.Cells.Item(11,7).Value = 0
$ws.Cells.Item(11,8).Value = 2.929109
$ws.Cells.Item(11,9).Value = 2.929109

# Row 12: INE522F01014 - Coal India Ltd
$ws.Cells.Item(12,1).Value = "INE522F01014"
$ws.Cells.Item(12,2).Value = "Coal India Ltd"
$ws.Cells.Item(12,3).Value = "quant Momentum Fund"
$ws.Cells.Item(12,4).Value = "Fresh Entry"
$ws.Cells.Item(12,5).Value = 2.211141
$ws.Cells.Item(12,6).Value = 0
$ws.Cells.Item(12,7).Value = 0
$ws.Cells.Item(12,8).Value = 2.211141
$ws.Cells.Item(12,9).Value = 2.211141

# Row 13: INE417T01026 - PB Fintech Limited
$ws.Cells.Item(13,1).Value = "INE417T01026"
$ws.Cells.Item(13,2).Value = "PB Fintech Limited"
$ws.Cells.Item(13,3).Value = "quant Momentum Fund"
$ws.Cells.Item(13,4).Value = "Fresh Entry"
$ws.Cells.Item(13,5).Value = 1.753087
$ws.Cells.Item(13,6).Value = 0
$ws.Cells.Item(13,7).Value = 0
$ws.Cells.Item(13,8).Value = 1.753087
$ws.Cells.Item(13,9).Value = 1.753087

# Row 14: INE075A01022 - Wipro Ltd
$ws.Cells.Item(14,1).Value = "INE075A01022"
$ws.Cells.Item(14,2).Value = "Wipro Ltd"
$ws.Cells.Item(14,3).Value = "quant Momentum Fund"
$ws.Cells.Item(14,4).Value = "Reducing"
$ws.Cells.Item(14,5).Value = 1.452308
$ws.Cells.Item(14,6).Value = 1.504776
$ws.Cells.Item(14,7).Value = 0
$ws.Cells.Item(14,8).Value = -0.05246799999999996
$ws.Cells.Item(14,9).Value = 1.452308

# Row 15: INE127D01025 - HDFC Asset Management Company Ltd
$ws.Cells.Item(15,1).Value = "INE127D01025"
$ws.Cells.Item(15,2).Value = "HDFC Asset Management Company Ltd"
$ws.Cells.Item(15,3).Value = "quant Momentum Fund"
$ws.Cells.Item(15,4).Value = "Fresh Entry"
$ws.Cells.Item(15,5).Value = 1.256642
$ws.Cells.Item(15,6).Value = 0
$ws.Cells.Item(15,7).Value = 0
$ws.Cells.Item(15,8).Value = 1.256642
$ws.Cells.Item(15,9).Value = 1.256642

# Row 16: INE259A01022 - Colgate-Palmolive (India) Ltd
$ws.Cells.Item(16,1).Value = "INE259A01022"
$ws.Cells.Item(16,2).Value = "Colgate-Palmolive (India) Ltd"
$ws.Cells.Item(16,3).Value = "quant Momentum Fund"
$ws.Cells.Item(16,4).Value = "Fresh Entry"
$ws.Cells.Item(16,5).Value = 0.897848
$ws.Cells.Item(16,6).Value = 0
$ws.Cells.Item(16,7).Value = 0
$ws.Cells.Item(16,8).Value = 0.897848
$ws.Cells.Item(16,9).Value = 0.897848

# Row 17: INE775A01035 - Samvardhana Motherson International Ltd
$ws.Cells.Item(17,1).Value = "INE775A01035"
$ws.Cells.Item(17,2).Value = "Samvardhana Motherson International Ltd"
$ws.Cells.Item(17,3).Value = "quant Momentum Fund"
$ws.Cells.Item(17,4).Value = "Reducing Consistently"
$ws.Cells.Item(17,5).Value = 0.360054
$ws.Cells.Item(17,6).Value = 7.212751
$ws.Cells.Item(17,7).Value = 3.452509
$ws.Cells.Item(17,8).Value = -6.852697
$ws.Cells.Item(17,9).Value = -3.092455

# Row 18: INE044A01036 - Sun Pharmaceutical Industries Limited
$ws.Cells.Item(18,1).Value = "INE044A01036"
$ws.Cells.Item(18,2).Value = "Sun Pharmaceutical Industries Limited"
$ws.Cells.Item(18,3).Value = "quant Momentum Fund"
$ws.Cells.Item(18,4).Value = "Complete Exit"
$ws.Cells.Item(18,5).Value = 0
$ws.Cells.Item(18,6).Value = 0
$ws.Cells.Item(18,7).Value = 3.422733
$ws.Cells.Item(18,8).Value = 0
$ws.Cells.Item(18,9).Value = -3.422733

# Row 19: INE647A01010 - SRF Limited
$ws.Cells.Item(19,1).Value = "INE647A01010"
$ws.Cells.Item(19,2).Value = "SRF Limited"
$ws.Cells.Item(19,3).Value = "quant Momentum Fund"
$ws.Cells.Item(19,4).Value = "Complete Exit"
$ws.Cells.Item(19,5).Value = 0
$ws.Cells.Item(19,6).Value = 0
$ws.Cells.Item(19,7).Value = 2.639061
$ws.Cells.Item(19,8).Value = 0
$ws.Cells.Item(19,9).Value = -2.639061

# Row 20: INE467B01029 - Tata Consultancy Services Limited
$ws.Cells.Item(20,1).Value = "INE467B01029"
$ws.Cells.Item(20,2).Value = "Tata Consultancy Services Limited"
$ws.Cells.Item(20,3).Value = "quant Momentum Fund"
$ws.Cells.Item(20,4).Value = "Complete Exit"
$ws.Cells.Item(20,5).Value = 0
$ws.Cells.Item(20,6).Value = 0.749743
$ws.Cells.Item(20,7).Value = 0
$ws.Cells.Item(20,8).Value = -0.749743
$ws.Cells.Item(20,9).Value = 0

# Row 21: INE424H01027 - SUN TV Network Limited
$ws.Cells.Item(21,1).Value = "INE424H01027"
$ws.Cells.Item(21,2).Value = "SUN TV Network Limited"
$ws.Cells.Item(21,3).Value = "quant Momentum Fund"
$ws.Cells.Item(21,4).Value = "Complete Exit"
$ws.Cells.Item(21,5).Value = 0
$ws.Cells.Item(21,6).Value = 0
$ws.Cells.Item(21,7).Value = 3.186287
$ws.Cells.Item(21,8).Value = 0
$ws.Cells.Item(21,9).Value = -3.186287

# Row 22: INE062A01020 - State Bank of India
$ws.Cells.Item(22,1).Value = "INE062A01020"
$ws.Cells.Item(22,2).Value = "State Bank of India"
$ws.Cells.Item(22,3).Value = "quant Momentum Fund"
$ws.Cells.Item(22,4).Value = "Complete Exit"
$ws.Cells.Item(22,5).Value = 0
$ws.Cells.Item(22,6).Value = 0
$ws.Cells.Item(22,7).Value = 3.063583
$ws.Cells.Item(22,8).Value = 0
$ws.Cells.Item(22,9).Value = -3.063583

# Row 23: INE094A01015 - Hindustan Petroleum Corporation Ltd
$ws.Cells.Item(23,1).Value = "INE094A01015"
$ws.Cells.Item(23,2).Value = "Hindustan Petroleum Corporation Ltd"
$ws.Cells.Item(23,3).Value = "quant Momentum Fund"
$ws.Cells.Item(23,4).Value = "Complete Exit"
$ws.Cells.Item(23,5).Value = 0
$ws.Cells.Item(23,6).Value = 2.976955
$ws.Cells.Item(23,7).Value = 0
$ws.Cells.Item(23,8).Value = -2.976955
$ws.Cells.Item(23,9).Value = 0

# Row 24: INE129A01019 - GAIL (India) Limited
$ws.Cells.Item(24,1).Value = "INE129A01019"
$ws.Cells.Item(24,2).Value = "GAIL (India) Limited"
$ws.Cells.Item(24,3).Value = "quant Momentum Fund"
$ws.Cells.Item(24,4).Value = "Complete Exit"
$ws.Cells.Item(24,5).Value = 0
$ws.Cells.Item(24,6).Value = 0
$ws.Cells.Item(24,7).Value = 1.719372
$ws.Cells.Item(24,8).Value = 0
$ws.Cells.Item(24,9).Value = -1.719372

# Row 25: INE280A01028 - Titan Company Limited
$ws.Cells.Item(25,1).Value = "INE280A01028"
$ws.Cells.Item(25,2).Value = "Titan Company Limited"
$ws.Cells.Item(25,3).Value = "quant Momentum Fund"
$ws.Cells.Item(25,4).Value = "Complete Exit"
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(25,6).Value = 2.93265
$ws.Cells.Item(25,7).Value = 0
$ws.Cells.Item(25,8).Value = -2.93265
$ws.Cells.Item(25,9).Value = 0

# Row 26: INE271C01023 - DLF Limited
$ws.Cells.Item(26,1).Value = "INE271C01023"
$ws.Cells.Item(26,2).Value = "DLF Limited"
$ws.Cells.Item(26,3).Value = "quant Momentum Fund"
$ws.Cells.Item(26,4).Value = "Complete Exit"
$ws.Cells.Item(26,5).Value = 0
$ws.Cells.Item(26,6).Value = 6.482575
$ws.Cells.Item(26,7).Value = 6.791959
$ws.Cells.Item(26,8).Value = -6.482575
$ws.Cells.Item(26,9).Value = -6.791959

# Row 27: INE245A01021 - Tata Power Company Limited
$ws.Cells.Item(27,1).Value = "INE245A01021"
$ws.Cells.Item(27,2).Value = "Tata Power Company Limited"
$ws.Cells.Item(27,3).Value = "quant Momentum Fund"
$ws.Cells.Item(27,4).Value = "Complete Exit"
$ws.Cells.Item(27,5).Value = 0
$ws.Cells.Item(27,6).Value = 3.680464
$ws.Cells.Item(27,7).Value = 7.396243
$ws.Cells.Item(27,8).Value = -3.680464
$ws.Cells.Item(27,9).Value = -7.396243

# Row 28: INE0CZ201020 - ANTHEM BIOSCIENCES LIMITED
$ws.Cells.Item(28,1).Value = "INE0CZ201020"
$ws.Cells.Item(28,2).Value = "ANTHEM BIOSCIENCES LIMITED"
$ws.Cells.Item(28,3).Value = "quant Momentum Fund"
$ws.Cells.Item(28,4).Value = "Complete Exit"
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(28,6).Value = 0
$ws.Cells.Item(28,7).Value = 6.446925
$ws.Cells.Item(28,8).Value = 0
$ws.Cells.Item(28,9).Value = -6.446925

# Row 29: INE237A01028 - Kotak Mahindra Bank Limited
$ws.Cells.Item(29,1).Value = "INE237A01028"
$ws.Cells.Item(29,2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(29,3).Value = "quant Momentum Fund"
$ws.Cells.Item(29,4).Value = "Complete Exit"
$ws.Cells.Item(29,5).Value = 0
$ws.Cells.Item(29,6).Value = 3.701086
$ws.Cells.Item(29,7).Value = 0
$ws.Cells.Item(29,8).Value = -3.701086
$ws.Cells.Item(29,9).Value = 0

# Row 30: INE0U4701011 - Digitide Solutions Limited
$ws.Cells.Item(30,1).Value = "INE0U4701011"
$ws.Cells.Item(30,2).Value = "Digitide Solutions Limited"
$ws.Cells.Item(30,3).Value = "quant Momentum Fund"
$ws.Cells.Item(30,4).Value = "Complete Exit"
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(30,6).Value = 0
$ws.Cells.Item(30,7).Value = 1.048161
$ws.Cells.Item(30,8).Value = 0
$ws.Cells.Item(30,9).Value = -1.048161

# Row 31: INE192A01025 - Tata Consumer Products Ltd
$ws.Cells.Item(31,1).Value = "INE192A01025"
$ws.Cells.Item(31,2).Value = "Tata Consumer Products Ltd"
$ws.Cells.Item(31,3).Value = "quant Momentum Fund"
$ws.Cells.Item(31,4).Value = "Complete Exit"
$ws.Cells.Item(31,5).Value = 0
$ws.Cells.Item(31,6).Value = 1.687128
$ws.Cells.Item(31,7).Value = 0
$ws.Cells.Item(31,8).Value = -1.687128
$ws.Cells.Item(31,9).Value = 0

# Row 32: INE030A01027 - Hindustan Unilever Limited
$ws.Cells.Item(32,1).Value = "INE030A01027"
$ws.Cells.Item(32,2).Value = "Hindustan Unilever Limited"
$ws.Cells.Item(32,3).Value = "quant Momentum Fund"
$ws.Cells.Item(32,4).Value = "Complete Exit"
$ws.Cells.Item(32,5).Value = 0
$ws.Cells.Item(32,6).Value = 1.46324
$ws.Cells.Item(32,7).Value = 0
$ws.Cells.Item(32,8).Value = -1.46324
$ws.Cells.Item(32,9).Value = 0
